$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Switch to manual calculation before editing so the engine's post-script
# auto-recalc does not recompute B6 (=SUM(C6+D6)) once D6 stops being a
# number below — the source workbook already carries a similarly "stale"
# cached formula result in B5 (whose D5 is the text "-4%"), so B6 should
# keep its existing cached value (935) rather than turn into #VALUE!.
$excel.Calculation = -4135

# Rename the "Total tea sales" column header (also updates the Table1
# column name bound to this header cell).
$ws.Range("B1").Value = "印度奶茶总销售数量（件）"

# Row 6 (2023-10-01): the handmade/pre-made tea sales figures were replaced
# with clock-time-looking text values.
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
